$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.937.45"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "3.190.78"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'537.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'144.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("D9").Value = "'7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +4.13%  "
$ws.Range("D11").Value = "'0.430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("D12").Value = "3.746.01"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "'26.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "59.999.38"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "3.188.82"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'13.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'8.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").Value = "'383.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").Value = "'70.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "'8.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.38%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "0.0₃0898"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'22.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("D31").Value = "'6.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "'5.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "'6.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("D35").Value = "'156.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").Value = "'1.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "2.775.00"
$ws.Range("E37").Value = "  +5.45%  "
$ws.Range("D38").Value = "'25.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "'4.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'39.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("D43").Value = "'0.729"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("D44").Value = "'0.0287"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.70%  "
$ws.Range("D45").Value = "3.234.36"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").Value = "'6.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "'0.801"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("D50").Value = "'20.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("E51").Value = "  +0.02%  "
